$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append after the existing last row (161)
$rows = @(
    @{ Row = 162; A = 12; B = "Mapocho Venta Directa de Santiago"; C = "Metropolitana"; D = 44628; E = 13; F = 100112043; G = "Pepino dulce"; H = "Cultivar IV Región"; I = "Especial"; J = 250; K = 14000; L = 14000; M = 14000; N = "`$/bandeja 18 kilos"; O = "Provincia de Limarí"; P = 778; Q = 18; R = "Hortaliza" },
    @{ Row = 163; A = 12; B = "Mapocho Venta Directa de Santiago"; C = "Metropolitana"; D = 44628; E = 13; F = 100112043; G = "Pepino dulce"; H = "Cultivar IV Región"; I = "Primera"; J = 220; K = 12000; L = 12000; M = 12000; N = "`$/bandeja 18 kilos"; O = "Provincia de Limarí"; P = 667; Q = 18; R = "Hortaliza" },
    @{ Row = 164; A = 12; B = "Mapocho Venta Directa de Santiago"; C = "Metropolitana"; D = 44628; E = 13; F = 100112043; G = "Pepino dulce"; H = "Cultivar IV Región"; I = "Segunda"; J = 280; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; O = "Provincia de Limarí"; P = 556; Q = 18; R = "Hortaliza" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C

    # Column D holds a date value, formatted like the preceding rows (format of D161)
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $r.D
    $dCell.NumberFormat = $ws.Cells.Item(161, 4).NumberFormat

    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
}
